$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30: 2021-02-17 (Excel serial 44244), all zero activity
$ws.Range("A30").Value = 44244
$ws.Range("B30").Value = 0
$ws.Range("C30").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 0

# Row 31: 2021-02-18 (Excel serial 44245)
$ws.Range("A31").Value = 44245
$ws.Range("B31").Value = 30
$ws.Range("C31").Value = 40
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 40
$ws.Range("F31").Value = 50

# Match the existing date formatting used by the rest of column A
$ws.Range("A29").Copy()
$ws.Range("A30:A31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection to reflect where the user ended up after entry
$ws.Range("G31").Select()
